# Add data for 2022-05-31
# Updates the workbook so the "through May 22" reporting date becomes
# "through May 23", and bumps the May-column carjacking counts for the
# neighborhoods that had an incident recorded on the newly-included day.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet and update the header label to reflect the new cutoff date.
$ws.Name = "Through 2022-05-23"
$ws.Range("B1").Value = "May 2022 (through May 23)"

# Bump existing counts by one additional carjacking each.
$ws.Range("G2").Value = 6    # Englewood, May 2021
$ws.Range("G5").Value = 4    # Garfield Park, May 2021
$ws.Range("AF5").Value = 8   # Garfield Park, May 2016
$ws.Range("AA6").Value = 2   # Chicago Lawn, May 2017
$ws.Range("L7").Value = 2    # North Lawndale, May 2020
$ws.Range("L11").Value = 3   # Roseland, May 2020
$ws.Range("G12").Value = 3   # Kenwood, May 2021
$ws.Range("AF25").Value = 3  # Auburn Gresham, May 2016
$ws.Range("Q27").Value = 2   # Wicker Park, May 2019
$ws.Range("B38").Value = 2   # Douglas, May 2022
$ws.Range("B45").Value = 3   # Logan Square, May 2022
$ws.Range("V46").Value = 2   # Little Village, May 2018
$ws.Range("B91").Value = 4   # Washington Park, May 2022

# Previously-empty cells that now have a first recorded carjacking.
$ws.Range("AF10").Value = 1  # Belmont Cragin, May 2016
$ws.Range("L16").Value = 1   # Little Italy, UIC, May 2020
$ws.Range("AA41").Value = 1  # Morgan Park, May 2017
$ws.Range("Q47").Value = 1   # Armour Square, May 2019
$ws.Range("G64").Value = 1   # Garfield Ridge, May 2021
